# Automatische test-sync: 2025-08-05 16:19:50
# Adds a new log entry (row 3) to the "Logs" sheet and a matching
# aggregate row (row 3) to the "Dashboard" sheet, then widens the
# conditional formatting ranges and the chart series references so
# they include the new row.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# ---------------------------------------------------------------
# 1. Append the new row of data to the "Logs" sheet (row 3)
# ---------------------------------------------------------------
$logs.Range("A3").Value = "Wil je deze klant bellen?"
$logs.Range("B3").Value = "mailmind.test@zohomail.eu"
$logs.Range("C3").Value = "Testmail #1: Wil je deze klant bellen?"
$logs.Range("D3").Value = "Klantenservice / Contact"
$logs.Range("E3").Value = "Bedankt, we hebben dit doorgestuurd naar klantenservice@bedrijf.nl."
$logs.Range("F3").Value = "2025-08-05 16:19:40"
$logs.Range("G3").Value = "Ja"
$logs.Range("H3").Value = "Ja"
$logs.Range("I3").Value = "Nee"
$logs.Range("J3").Value = "Nee"

# ---------------------------------------------------------------
# 2. Append the matching aggregate row to the "Dashboard" sheet (row 3)
# ---------------------------------------------------------------
$dash.Range("A3").Value = "Klantenservice / Contact"
$dash.Range("B3").Value = 1

# ---------------------------------------------------------------
# 3. Widen the conditional formatting ranges on "Logs" so that each
#    rule also covers the newly added row 3.
# ---------------------------------------------------------------
$colRanges = @("D2", "G2", "H2", "I2", "J2")
foreach ($col in $colRanges) {
    $newRange = $logs.Range($col + ":" + ($col -replace "2","3"))
    $fcs = $logs.Range($col).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        $fc.ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------
# 4. Update the Dashboard bar chart so its category/value series
#    reference the widened data range (A2:A3 / B2:B3).
# ---------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$3,Dashboard!`$B`$2:`$B`$3,1)"
